# pnl, bs, cf are generated in excel through /routes
#
# The P&L sheet stores its figures as literal text strings (e.g. "330,858.46"),
# not numbers, so the updated figures below must also land back in the cells
# as text (not get auto-parsed into numeric values by Excel's smart-typing).
#
# Plain `Range.Value = "7,677,623.38"` would be interpreted by Excel as a
# number (stripping the thousands separators / trailing zeros and changing
# the cell's number format), so each write:
#   1. Forces the target cell to Text format ("@") before assigning the
#      value, so the literal string is preserved verbatim.
#   2. Re-applies the cell's original (General) formatting afterwards by
#      copying formats-only from an untouched neighboring cell that still
#      carries the original style, so the cell's style index is unaffected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Text,
        [string]$FormatDonor
    )
    $ws.Range($Address).NumberFormat = "@"
    $ws.Range($Address).Value = $Text
    $ws.Range($FormatDonor).Copy()
    $ws.Range($Address).PasteSpecial(-4122)  # xlPasteFormats
}

# Row 4 - I. Revenue from Operations
Set-TextValue "C4" "7,677,623.38" "B3"

# Row 5 - II. Other Income
Set-TextValue "C5" "0.00" "B3"

# Row 6 - III. Total Income (I + II)
Set-TextValue "C6" "7,677,623.38" "B3"

# Row 8 - Cost of Materials Consumed
Set-TextValue "C8" "2,498.09" "B8"

# Row 9 - Employee Benefit Expense
Set-TextValue "C9" "63,498,614.00" "B8"

# Row 10 - Other Expenses
Set-TextValue "C10" "215,212,520.99" "B8"

# Row 11 - Depreciation and Amortisation Expense (values swapped between C/D)
Set-TextValue "C11" "0.00" "B8"
Set-TextValue "D11" "60,344,362.00" "B8"

# Row 13 - Finance Costs
Set-TextValue "C13" "10,242,315.65" "B8"

# Row 15 - Total Expenses
Set-TextValue "C15" "288,955,948.73" "B8"
Set-TextValue "D15" "60,344,362.00" "B8"

# Row 16 - V. Profit Before Exceptional and Extraordinary Items and Tax (III - IV)
Set-TextValue "C16" "-281,278,325.35" "B3"
Set-TextValue "D16" "-60,344,362.00" "B3"

# Row 18 - VII. Profit Before Tax (V - VI)
Set-TextValue "C18" "-281,278,325.35" "B3"
Set-TextValue "D18" "-60,344,362.00" "B3"

# Row 22 - IX. Profit After Tax for the period (VII - VIII)
Set-TextValue "C22" "-281,278,325.35" "B3"
Set-TextValue "D22" "-60,344,362.00" "B3"
